$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "301.91"
$ws.Range("E2").Value = "'1.49%"
$ws.Range("D3").Value = "32.25"
$ws.Range("E3").Value = "'2.55%"
$ws.Range("D4").Value = "5.134"
$ws.Range("D5").Value = "0.07823"
$ws.Range("E5").Value = "'-2.23%"
$ws.Range("D6").Value = "2.263"
$ws.Range("E6").Value = "'-7.83%"
$ws.Range("D7").Value = "7.813"
$ws.Range("E7").Value = "'0.30%"
$ws.Range("D8").Value = "'3.810"
$ws.Range("E8").Value = "'0.36%"
$ws.Range("D9").Value = "0.9296"
$ws.Range("E9").Value = "'0.75%"
$ws.Range("D10").Value = "0.1758"
$ws.Range("E10").Value = "'-0.02%"
$ws.Range("D11").Value = "0.07684"
$ws.Range("E11").Value = "'4.79%"
$ws.Range("D12").Value = "0.08811"
$ws.Range("E12").Value = "'-0.02%"
$ws.Range("D13").Value = "'0.03060"
$ws.Range("E13").Value = "'0.82%"
$ws.Range("D14").Value = "0.1001"
$ws.Range("E14").Value = "'0.11%"
$ws.Range("D15").Value = "0.001526"
$ws.Range("E15").Value = "'1.53%"
$ws.Range("D16").Value = "0.006021"
$ws.Range("E16").Value = "'1.64%"
$ws.Range("D17").Value = "3.464"
$ws.Range("E17").Value = "'-1.08%"
$ws.Range("D18").Value = "2.254"
$ws.Range("E18").Value = "'0.27%"
$ws.Range("D20").Value = "0.1347"
$ws.Range("E20").Value = "'0.80%"
$ws.Range("E21").Value = "'-0.16%"
$ws.Range("D22").Value = "0.1822"
$ws.Range("E22").Value = "'12.73%"
$ws.Range("D23").Value = "0.04618"
$ws.Range("E23").Value = "'0.45%"
$ws.Range("D24").Value = "0.001255"
$ws.Range("E24").Value = "'0.98%"
$ws.Range("D25").Value = "0.004508"
$ws.Range("E25").Value = "'1.71%"
$ws.Range("D26").Value = "0.0001254"
$ws.Range("E26").Value = "'4.45%"
$ws.Range("D39").Value = "0.01791"
$ws.Range("E39").Value = "'1.00%"
$ws.Range("D40").Value = "0.04699"
$ws.Range("E40").Value = "'5.29%"
$ws.Range("D41").Value = "0.007243"
$ws.Range("E41").Value = "'3.96%"
$ws.Range("D42").Value = "0.1376"
$ws.Range("E42").Value = "'2.36%"
$ws.Range("D43").Value = "0.002128"
$ws.Range("E43").Value = "'-3.75%"
$ws.Range("D44").Value = "0.01109"
$ws.Range("E44").Value = "'13.11%"
$ws.Range("D45").Value = "0.00006339"
$ws.Range("E45").Value = "'-3.38%"
$ws.Range("D46").Value = "0.00000000752"
$ws.Range("E46").Value = "'0.15%"
$ws.Range("D47").Value = "0.003207"
$ws.Range("E47").Value = "'-38.76%"
$ws.Range("D48").Value = "0.7477"
$ws.Range("E48").Value = "'-8.88%"
$ws.Range("E49").Value = "'0.15%"
$ws.Range("E50").Value = "'0.15%"
